$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing rows 41-45 with revised figures ---

# Row 41
$ws.Range("B41").Value = 2169902
$ws.Range("C41").Value = 971223
$ws.Range("F41").Value = 815372
$ws.Range("G41").Value = 8149
$ws.Range("I41").Value = 101006
$ws.Range("J41").Value = 1789742
$ws.Range("K41").Value = 1070606
$ws.Range("L41").Value = 413519
$ws.Range("N41").Value = 261374
$ws.Range("Q41").Value = 380161
$ws.Range("R41").Value = 100110
$ws.Range("T41").Value = 103848
$ws.Range("W41").Value = 2177568
$ws.Range("X41").Value = 1897517
$ws.Range("Y41").Value = 280050

# Row 42
$ws.Range("B42").Value = 1875813
$ws.Range("C42").Value = 661212
$ws.Range("F42").Value = 762332
$ws.Range("I42").Value = 75872
$ws.Range("J42").Value = 1835746
$ws.Range("K42").Value = 1103012
$ws.Range("L42").Value = 441599
$ws.Range("N42").Value = 269368
$ws.Range("P42").Value = 10094
$ws.Range("Q42").Value = 40067
$ws.Range("R42").Value = 106969
$ws.Range("T42").Value = 104476
$ws.Range("W42").Value = 1876118
$ws.Range("X42").Value = 1943020
$ws.Range("Y42").Value = -66902

# Row 43
$ws.Range("B43").Value = 2052091
$ws.Range("C43").Value = 912736
$ws.Range("F43").Value = 745986
$ws.Range("G43").Value = 6656
$ws.Range("I43").Value = 98217
$ws.Range("J43").Value = 1795664
$ws.Range("K43").Value = 1109490
$ws.Range("L43").Value = 433335
$ws.Range("N43").Value = 229302
$ws.Range("P43").Value = 12719
$ws.Range("Q43").Value = 256428
$ws.Range("W43").Value = 2054669
$ws.Range("X43").Value = 1899827
$ws.Range("Y43").Value = 154842

# Row 44
$ws.Range("B44").Value = 2147093
$ws.Range("C44").Value = 849517
$ws.Range("F44").Value = 827462
$ws.Range("H44").Value = 336853
$ws.Range("I44").Value = 125700
$ws.Range("J44").Value = 1996637
$ws.Range("K44").Value = 1182728
$ws.Range("L44").Value = 520871
$ws.Range("M44").Value = 392
$ws.Range("N44").Value = 256151
$ws.Range("O44").Value = 18568
$ws.Range("P44").Value = 17928
$ws.Range("Q44").Value = 150456
$ws.Range("R44").Value = 150938
$ws.Range("T44").Value = 149190
$ws.Range("U44").Value = 4372
$ws.Range("W44").Value = 2149716
$ws.Range("X44").Value = 2150199
$ws.Range("Y44").Value = -483

# Row 45
$ws.Range("B45").Value = 2324173
$ws.Range("C45").Value = 1074335
$ws.Range("F45").Value = 819926
$ws.Range("G45").Value = 6255
$ws.Range("H45").Value = 312139
$ws.Range("I45").Value = 111517
$ws.Range("J45").Value = 1804283
$ws.Range("K45").Value = 1103040
$ws.Range("L45").Value = 387928
$ws.Range("M45").Value = 317
$ws.Range("N45").Value = 266870
$ws.Range("O45").Value = 28892
$ws.Range("P45").Value = 17235
$ws.Range("Q45").Value = 519890
$ws.Range("R45").Value = 106861
$ws.Range("S45").Value = 1202
$ws.Range("T45").Value = 104672
$ws.Range("U45").Value = 3390
$ws.Range("W45").Value = 2325375
$ws.Range("X45").Value = 1912345
$ws.Range("Y45").Value = 413030

# --- Append new row 46 (01-04-2021) ---
# Assign via a text formula first (so Excel stores it as a text string,
# not auto-converted to a date serial), then paste-special as values only
# so the formula collapses to a plain shared-string cell with no extra
# number formatting/style applied.
$ws.Range("A46").Formula = '="01-04-2021"'
$ws.Range("A46").Copy()
$ws.Range("A46").PasteSpecial(-4163)
$ws.Range("B46").Value = 2056989
$ws.Range("C46").Value = 739748
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 795986
$ws.Range("G46").Value = 6364
$ws.Range("H46").Value = 400930
$ws.Range("I46").Value = 113961
$ws.Range("J46").Value = 1988083
$ws.Range("K46").Value = 1165582
$ws.Range("L46").Value = 474178
$ws.Range("M46").Value = 382
$ws.Range("N46").Value = 311035
$ws.Range("O46").Value = 22373
$ws.Range("P46").Value = 14533
$ws.Range("Q46").Value = 68906
$ws.Range("R46").Value = 136842
$ws.Range("S46").Value = 1030
$ws.Range("T46").Value = 135172
$ws.Range("U46").Value = 2701
$ws.Range("W46").Value = 2058020
$ws.Range("X46").Value = 2125955
$ws.Range("Y46").Value = -67936
